$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.915.34'
$ws.Range('E2').Value = '  -0.38%  '
$ws.Range('D3').Value = '3.862.01'
$ws.Range('E3').Value = '  +1.07%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'469.52"
$ws.Range('E5').Value = '  +3.74%  '
$ws.Range('D6').Value = "'144.61"
$ws.Range('E6').Value = '  -1.07%  '
$ws.Range('D7').Value = "'0.610"
$ws.Range('E7').Value = '  -1.77%  '
$ws.Range('D8').Value = "'0.998"
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('E9').Value = '  -3.94%  '
$ws.Range('E10').Value = '  +5.80%  '
$ws.Range('D11').Value = "'0.0000345"
$ws.Range('E11').Value = '  +8.44%  '
$ws.Range('D12').Value = "'41.89"
$ws.Range('E12').Value = '  -4.10%  '
$ws.Range('D13').Value = '4.475.93'
$ws.Range('E13').Value = '  +0.87%  '
$ws.Range('E14').Value = '  -2.36%  '
$ws.Range('D15').Value = '3.927.83'
$ws.Range('E15').Value = '  +3.36%  '
$ws.Range('E16').Value = '  -2.92%  '
$ws.Range('E17').Value = '  -0.47%  '
$ws.Range('D18').Value = "'19.54"
$ws.Range('E18').Value = '  -2.38%  '
$ws.Range('E19').Value = '  -4.76%  '
$ws.Range('D20').Value = '66.955.58'
$ws.Range('E20').Value = '  -0.48%  '
$ws.Range('D21').Value = "'426.18"
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('E22').Value = '  +1.32%  '
$ws.Range('D23').Value = "'14.16"
$ws.Range('E23').Value = '  -4.47%  '
$ws.Range('D24').Value = "'85.96"
$ws.Range('E24').Value = '  -0.43%  '
$ws.Range('E25').Value = '  +0.60%  '
$ws.Range('D26').Value = "'37.55"
$ws.Range('E26').Value = '  +0.99%  '
$ws.Range('D27').Value = "'9.99"
$ws.Range('E27').Value = '  -2.95%  '
$ws.Range('D28').Value = "'10.07"
$ws.Range('E28').Value = '  +3.43%  '
$ws.Range('D29').Value = "'720.95"
$ws.Range('E29').Value = '  -2.11%  '
$ws.Range('D30').Value = "'13.09"
$ws.Range('E30').Value = '  -5.13%  '
$ws.Range('E31').Value = '  -5.36%  '
$ws.Range('E32').Value = '  +2.19%  '
$ws.Range('D33').Value = "'41.55"
$ws.Range('E33').Value = '  -3.43%  '
$ws.Range('D34').Value = '0.0₃0870'
$ws.Range('E34').Value = '  +26.37%  '
$ws.Range('D35').Value = "'58.21"
$ws.Range('E35').Value = '  +1.78%  '
$ws.Range('D36').Value = "'0.150"
$ws.Range('E36').Value = '  -5.18%  '
$ws.Range('D37').Value = "'0.998"
$ws.Range('D38').Value = "'5.29"
$ws.Range('E38').Value = '  -4.99%  '
$ws.Range('E39').Value = '  -2.84%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').Value = "'2.73"
$ws.Range('E40').Value = '  +4.90%  '
$ws.Range('B41').Value = 'ThetaToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D41').Value = "'2.99"
$ws.Range('E41').Value = '  +2.73%  '
$ws.Range('D42').Value = "'2.97"
$ws.Range('E42').Value = '  +10.83%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = "'1.00"
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').Value = "'0.336"
$ws.Range('E44').Value = '  -3.53%  '
$ws.Range('E45').Value = '  -1.14%  '
$ws.Range('D46').Value = "'3.40"
$ws.Range('E46').Value = '  -1.44%  '
$ws.Range('D47').Value = "'2.14"
$ws.Range('E47').Value = '  +0.25%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').Value = "'144.79"
$ws.Range('E48').Value = '  +0.36%  '
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').Value = "'3.12"
$ws.Range('E49').Value = '  -4.20%  '
$ws.Range('D50').Value = "'2.79"
$ws.Range('E50').Value = '  -2.47%  '
$ws.Range('D51').Value = "'24.29"
$ws.Range('E51').Value = '  -2.74%  '
